$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1149921.2
$ws.Range("I11").Value = 1149921.2
$ws.Range("K11").Value = 1149921.2
$ws.Range("M11").Value = -1149781.2
$ws.Range("H15").Value = 226.09
$ws.Range("I15").Value = 226.09
$ws.Range("K15").Value = 678.27
$ws.Range("M15").Value = -509.27
$ws.Range("H28").Value = 583.88
$ws.Range("I28").Value = 562.6
$ws.Range("J28").Value = 669
$ws.Range("K28").Value = 562.6
$ws.Range("L28").Value = 669
$ws.Range("M28").Value = -77.60000000000002
$ws.Range("N28").Value = -1639
$ws.Range("H44").Value = 22080.77
$ws.Range("J44").Value = 22080.77
$ws.Range("L44").Value = 22080.77
$ws.Range("N44").Value = -23004.77
$ws.Range("H58").Value = 8809.154
$ws.Range("I58").Value = 1614.3334
$ws.Range("J58").Value = 24997.5
$ws.Range("K58").Value = 4843.0002
$ws.Range("L58").Value = 74992.5
$ws.Range("M58").Value = -4693.0002
$ws.Range("N58").Value = -75292.5
$ws.Range("H87").Value = 23318.5
$ws.Range("J87").Value = 23318.5
$ws.Range("L87").Value = 23318.5
$ws.Range("N87").Value = -25814.5
$ws.Range("H90").Value = 23318.5
$ws.Range("J90").Value = 23318.5
$ws.Range("L90").Value = 69955.5
$ws.Range("N90").Value = -82435.5
$ws.Range("H100").Value = 18183218
$ws.Range("I100").Value = 20001350
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 20001350
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -20000809
$ws.Range("N100").Value = -2982
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H141").Value = 4897.1562
$ws.Range("I141").Value = 4953.893
$ws.Range("K141").Value = 14861.679
$ws.Range("M141").Value = -9681.679

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6276.1826
$ws.Range("I32").Value = 3780.0166
$ws.Range("J32").Value = 10814.667
$ws.Range("K32").Value = 3780.0166
$ws.Range("L32").Value = 10814.667
$ws.Range("M32").Value = -3493.0166
$ws.Range("N32").Value = -11388.667
$ws.Range("H74").Value = 4138.3228
$ws.Range("I74").Value = 4237.3335
$ws.Range("J74").Value = 3798.8572
$ws.Range("K74").Value = 4237.3335
$ws.Range("L74").Value = 3798.8572
$ws.Range("M74").Value = -3363.3335
$ws.Range("N74").Value = -5546.8572
$ws.Range("H77").Value = 4138.3228
$ws.Range("I77").Value = 4237.3335
$ws.Range("J77").Value = 3798.8572
$ws.Range("K77").Value = 21186.6675
$ws.Range("L77").Value = 18994.286
$ws.Range("M77").Value = -16818.6675
$ws.Range("N77").Value = -27730.286
$ws.Range("H137").Value = 52770
$ws.Range("J137").Value = 52770
$ws.Range("L137").Value = 52770
$ws.Range("N137").Value = -62970

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 18500
$ws.Range("J6").Value = 18500
$ws.Range("L6").Value = 18500
$ws.Range("N6").Value = -18726
$ws.Range("H134").Value = 3245.6743
$ws.Range("I134").Value = 1712.1875
$ws.Range("J134").Value = 7706.727
$ws.Range("K134").Value = 5136.5625
$ws.Range("L134").Value = 23120.181
$ws.Range("M134").Value = -2601.5625
$ws.Range("N134").Value = -28190.181
$ws.Range("H137").Value = 32948
$ws.Range("J137").Value = 32948
$ws.Range("L137").Value = 32948
$ws.Range("N137").Value = -43148
$ws.Range("H139").Value = 40780
$ws.Range("J139").Value = 40780
$ws.Range("L139").Value = 40780
$ws.Range("N139").Value = -51060

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45335
$ws.Range("J20").Value = 45335
$ws.Range("L20").Value = 45335
$ws.Range("N20").Value = -45807
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H30").Value = 45335
$ws.Range("J30").Value = 45335
$ws.Range("L30").Value = 45335
$ws.Range("N30").Value = -45517
$ws.Range("H31").Value = 3018.0652
$ws.Range("I31").Value = 1274.6957
$ws.Range("J31").Value = 4761.4346
$ws.Range("K31").Value = 1274.6957
$ws.Range("L31").Value = 4761.4346
$ws.Range("M31").Value = -979.6957
$ws.Range("N31").Value = -5351.4346
$ws.Range("H34").Value = 3018.0652
$ws.Range("I34").Value = 1274.6957
$ws.Range("J34").Value = 4761.4346
$ws.Range("K34").Value = 1274.6957
$ws.Range("L34").Value = 4761.4346
$ws.Range("M34").Value = -1072.6957
$ws.Range("N34").Value = -5165.4346
$ws.Range("H36").Value = 16809.6
$ws.Range("I36").Value = 8016
$ws.Range("J36").Value = 30000
$ws.Range("K36").Value = 8016
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = -7628
$ws.Range("N36").Value = -30776
$ws.Range("H40").Value = 16809.6
$ws.Range("I40").Value = 8016
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = 8016
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = -7856
$ws.Range("N40").Value = -30320
$ws.Range("H120").Value = 29999.076
$ws.Range("J120").Value = 29999.076
$ws.Range("L120").Value = 29999.076
$ws.Range("N120").Value = -37257.076
$ws.Range("H123").Value = 36374.547
$ws.Range("J123").Value = 36374.547
$ws.Range("L123").Value = 36374.547
$ws.Range("N123").Value = -46174.547
$ws.Range("H128").Value = 45335
$ws.Range("J128").Value = 45335
$ws.Range("L128").Value = 45335
$ws.Range("N128").Value = -55295
$ws.Range("H132").Value = 2796.2222
$ws.Range("I132").Value = 1357.6666
$ws.Range("J132").Value = 4234.778
$ws.Range("K132").Value = 4072.9998
$ws.Range("L132").Value = 12704.334
$ws.Range("M132").Value = -1542.9998
$ws.Range("N132").Value = -17764.334
$ws.Range("H134").Value = 5245.033
$ws.Range("I134").Value = 6702.778
$ws.Range("J134").Value = 3058.4167
$ws.Range("K134").Value = 20108.334
$ws.Range("L134").Value = 9175.250100000001
$ws.Range("M134").Value = -17573.334
$ws.Range("N134").Value = -14245.2501

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 66.07143000000001
$ws.Range("I12").Value = 173.33333
$ws.Range("J12").Value = 36.81818
$ws.Range("K12").Value = 519.99999
$ws.Range("L12").Value = 110.45454
$ws.Range("M12").Value = -346.99999
$ws.Range("N12").Value = -456.45454
$ws.Range("H113").Value = 727.4651
$ws.Range("I113").Value = 638.2
$ws.Range("K113").Value = 1914.6
$ws.Range("M113").Value = 255.3999999999999
$ws.Range("H131").Value = 11364639
$ws.Range("I131").Value = 100002260
$ws.Range("J131").Value = 841.5641000000001
$ws.Range("K131").Value = 300006780
$ws.Range("L131").Value = 2524.6923
$ws.Range("M131").Value = -300001740
$ws.Range("N131").Value = -12604.6923

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20490
$ws.Range("H99").Value = 4333.3335
$ws.Range("I99").Value = 4333.3335
$ws.Range("K99").Value = 4333.3335
$ws.Range("M99").Value = -2087.3335
$ws.Range("H137").Value = 82526.664
$ws.Range("J137").Value = 82526.664
$ws.Range("L137").Value = 82526.664
$ws.Range("N137").Value = -92726.664

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4380.0557
$ws.Range("I7").Value = 2094.2
$ws.Range("J7").Value = 7237.375
$ws.Range("K7").Value = 2094.2
$ws.Range("L7").Value = 7237.375
$ws.Range("M7").Value = -1982.2
$ws.Range("N7").Value = -7461.375
$ws.Range("H93").Value = 2630.7222
$ws.Range("I93").Value = 1946.0834
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 1946.0834
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -698.0834
$ws.Range("N93").Value = -6496
$ws.Range("H122").Value = 3189.4443
$ws.Range("I122").Value = 2871.2058
$ws.Range("K122").Value = 8613.617400000001
$ws.Range("M122").Value = -6163.617400000001
$ws.Range("H126").Value = 4380.0557
$ws.Range("I126").Value = 2094.2
$ws.Range("J126").Value = 7237.375
$ws.Range("K126").Value = 6282.599999999999
$ws.Range("L126").Value = 21712.125
$ws.Range("M126").Value = -3812.599999999999
$ws.Range("N126").Value = -26652.125
$ws.Range("H140").Value = 67714.11
$ws.Range("J140").Value = 67714.11
$ws.Range("L140").Value = 67714.11
$ws.Range("N140").Value = -78074.11

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17858642
$ws.Range("I81").Value = 23810856
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 47621712
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -47620651
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 17858642
$ws.Range("I84").Value = 23810856
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 238108560
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -238103256
$ws.Range("N84").Value = -30608
$ws.Range("H126").Value = 282414.6
$ws.Range("I126").Value = 1810.0435
$ws.Range("J126").Value = 712674.9399999999
$ws.Range("K126").Value = 5430.1305
$ws.Range("L126").Value = 2138024.82
$ws.Range("M126").Value = -2960.1305
$ws.Range("N126").Value = -2142964.82
$ws.Range("H130").Value = 39116
$ws.Range("J130").Value = 39116
$ws.Range("L130").Value = 39116
$ws.Range("N130").Value = -49156
$ws.Range("H136").Value = 3973.08
$ws.Range("I136").Value = 794.64703
$ws.Range("J136").Value = 10727.25
$ws.Range("K136").Value = 2383.94109
$ws.Range("L136").Value = 32181.75
$ws.Range("M136").Value = 166.0589100000002
$ws.Range("N136").Value = -37281.75
